# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The sheet's column G (header "K") held stale "strike#" counts. This
# regenerates that column with the freshly computed K values for each
# data row (rows 2-73), leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(1,2,2,1,1,0,2,1,0,1,2,1,1,1,0,1,0,2,1,1,1,0,0,1,2,2,3,1,3,2,1,0,1,0,3,0,1,0,1,1,2,0,2,0,0,1,0,1,0,1,0,0,1,2,1,1,0,0,1,2,1,2,0,1,1,1,0,4,1,1,1,0)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
